# Updates the Guca2a-Gucy2c LR-pairs sheet with newly-computed TPM values.
# The "Target cluster" dimension gains "Inflammatory-Mac" as a third option
# (in addition to "FAPs" and "Resolving-Mac") for every "Sending cluster",
# so rows 2-9 get new values and rows 10-13 are newly added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Guca2a", "Gucy2c", "FAPs", 3, 1, 1.134962, 3.404886, 0.4107584930630159, 0.4107584930630159, 2, 0.6666666666666666, 0.1264736666666667, 0.379421, 0.1713830533493957, 0.1713830533493957, 0.1435428056673334, 1.291885251006, 0.07039704473033621, 0.07039704473033621),
    @("ECs", "Guca2a", "Gucy2c", "Inflammatory-Mac", 3, 1, 1.134962, 3.404886, 0.4107584930630159, 0.4107584930630159, 2, 0.6666666666666666, 0.5694333333333333, 1.7083, 0.7716327510516618, 0.7716327510516618, 0.6462851948666667, 5.8165667538, 0.3169547060200499, 0.3169547060200499),
    @("ECs", "Guca2a", "Gucy2c", "Resolving-Mac", 3, 1, 1.134962, 3.404886, 0.4107584930630159, 0.4107584930630159, 1, 0.3333333333333333, 0.042052, 0.126156, 0.05698419559894249, 0.05698419559894248, 0.047727422024, 0.429546798216, 0.02340674231262976, 0.02340674231262975),
    @("FAPs", "Guca2a", "Gucy2c", "FAPs", 2, 0.6666666666666666, 0.5207596666666666, 1.562279, 0.1884701478357852, 0.1884701478357852, 2, 0.6666666666666666, 0.1264736666666667, 0.379421, 0.1713830533493957, 0.1713830533493957, 0.06586238449544445, 0.592761460459, 0.03230058940130886, 0.03230058940130886),
    @("FAPs", "Guca2a", "Gucy2c", "Inflammatory-Mac", 2, 0.6666666666666666, 0.5207596666666666, 1.562279, 0.1884701478357852, 0.1884701478357852, 2, 0.6666666666666666, 0.5694333333333333, 1.7083, 0.7716327510516618, 0.7716327510516618, 0.2965379128555555, 2.6688412157, 0.1454297386656404, 0.1454297386656404),
    @("FAPs", "Guca2a", "Gucy2c", "Resolving-Mac", 2, 0.6666666666666666, 0.5207596666666666, 1.562279, 0.1884701478357852, 0.1884701478357852, 1, 0.3333333333333333, 0.042052, 0.126156, 0.05698419559894249, 0.05698419559894248, 0.02189898550266666, 0.197090869524, 0.01073981976883599, 0.01073981976883599),
    @("Inflammatory-Mac", "Guca2a", "Gucy2c", "FAPs", 3, 1, 0.8996876666666666, 2.699063, 0.3256094478822911, 0.3256094478822911, 2, 0.6666666666666666, 0.1264736666666667, 0.379421, 0.1713830533493957, 0.1713830533493957, 0.1137867980581111, 1.024081182523, 0.05580394137747797, 0.05580394137747795),
    @("Inflammatory-Mac", "Guca2a", "Gucy2c", "Inflammatory-Mac", 3, 1, 0.8996876666666666, 2.699063, 0.3256094478822911, 0.3256094478822911, 2, 0.6666666666666666, 0.5694333333333333, 1.7083, 0.7716327510516618, 0.7716327510516618, 0.5123121469888888, 4.6108093229, 0.251250914037825, 0.251250914037825),
    @("Inflammatory-Mac", "Guca2a", "Gucy2c", "Resolving-Mac", 3, 1, 0.8996876666666666, 2.699063, 0.3256094478822911, 0.3256094478822911, 1, 0.3333333333333333, 0.042052, 0.126156, 0.05698419559894249, 0.05698419559894248, 0.03783366575866666, 0.3405029918279999, 0.01855459246698815, 0.01855459246698814),
    @("MuSCs", "Guca2a", "Gucy2c", "FAPs", 2, 0.6666666666666666, 0.207679, 0.6230370000000001, 0.07516191121890785, 0.07516191121890783, 2, 0.6666666666666666, 0.1264736666666667, 0.379421, 0.1713830533493957, 0.1713830533493957, 0.02626592461966667, 0.236393321577, 0.01288147784027262, 0.01288147784027262),
    @("MuSCs", "Guca2a", "Gucy2c", "Inflammatory-Mac", 2, 0.6666666666666666, 0.207679, 0.6230370000000001, 0.07516191121890785, 0.07516191121890783, 2, 0.6666666666666666, 0.5694333333333333, 1.7083, 0.7716327510516618, 0.7716327510516618, 0.1182593452333334, 1.0643341071, 0.05799739232814663, 0.05799739232814662),
    @("MuSCs", "Guca2a", "Gucy2c", "Resolving-Mac", 2, 0.6666666666666666, 0.207679, 0.6230370000000001, 0.07516191121890785, 0.07516191121890783, 1, 0.3333333333333333, 0.042052, 0.126156, 0.05698419559894249, 0.05698419559894248, 0.008733317308000001, 0.078599855772, 0.004283041050488595, 0.004283041050488593)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $rowValues[$c]
    }
}

Write-Output "Updated $($data.Count) rows starting at row $startRow"
